# Auto-generated Excel COM-interop script
# Applies the "cryptos list" GitHub Actions price/volume refresh: updates
# Price (D) and Volume(1h) (E) text values across rows 2-51, and restores
# the correct row order for "Dai" (row 29) / "Fetch.AI" (row 30), which had
# been swapped.
#
# All D/E columns hold plain text (inlineStr) in the source workbook, not
# numbers. Values that look numeric (e.g. "706.89", "1.00") are written
# with a leading apostrophe so Excel stores them as literal text instead
# of silently coercing them to numbers (which would also destroy
# significant trailing/leading zeros such as "1.00" or "0.730").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range('D2').Value = '70.922.97'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '3.821.27'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''706.89'
$ws.Range('E5').Value = '  +2.57%  '
$ws.Range('D6').Value = '''170.25'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('D7').Value = '3.821.32'
$ws.Range('E7').Value = '  -0.60%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('D11').Value = '''7.37'
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').Value = '''0.456'
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('D13').Value = '''0.0000253'
$ws.Range('E13').Value = '  -1.73%  '
$ws.Range('D14').Value = '''36.39'
$ws.Range('E14').Value = '  -0.58%  '
$ws.Range('D15').Value = '4.467.29'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '3.938.48'
$ws.Range('E16').Value = '  +2.59%  '
$ws.Range('D17').Value = '70.988.32'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').Value = '''17.30'
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('D21').Value = '''494.24'
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('D22').Value = '''10.60'
$ws.Range('E22').Value = '  -4.58%  '
$ws.Range('D23').Value = '''0.730'
$ws.Range('E23').Value = '  +1.45%  '
$ws.Range('D24').Value = '''85.66'
$ws.Range('E24').Value = '  +1.20%  '
$ws.Range('E25').Value = '  -1.78%  '
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').Value = '''12.07'
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('D28').Value = '3.975.81'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = '''2.07'
$ws.Range('E30').Value = '  -3.53%  '
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('D32').Value = '''7.38'
$ws.Range('E32').Value = '  -3.19%  '
$ws.Range('E33').Value = '  -3.60%  '
$ws.Range('D34').Value = '''29.27'
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('D35').Value = '''0.173'
$ws.Range('E35').Value = '  -4.86%  '
$ws.Range('D36').Value = '3.791.05'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').Value = '''9.12'
$ws.Range('E37').Value = '  -1.57%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').Value = '  -1.87%  '
$ws.Range('E40').Value = '  +3.47%  '
$ws.Range('E41').Value = '  -3.02%  '
$ws.Range('D42').Value = '''5.95'
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('D43').Value = '''3.28'
$ws.Range('E43').Value = '  -4.06%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = '''163.80'
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('E47').Value = '  +1.23%  '
$ws.Range('D48').Value = '''427.56'
$ws.Range('E48').Value = '  +3.23%  '
$ws.Range('D49').Value = '''48.82'
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('E51').Value = '  -2.62%  '
